$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12329.2895096546
$ws.Range("C2").Value = 11750.7392695781
$ws.Range("E2").Value = 8121.63873300613
$ws.Range("F2").Value = -5.14508322565659

$ws.Range("B3").Value = 12602.479458358
$ws.Range("C3").Value = 12161.631818983
$ws.Range("E3").Value = 8193.75913948673
$ws.Range("F3").Value = 344.980456602904

$ws.Range("B4").Value = 12495.4648953028
$ws.Range("C4").Value = 11701.7312803431
$ws.Range("E4").Value = 8401.12438578046
$ws.Range("F4").Value = 334.458152755149

$ws.Range("B5").Value = 12406.8406380375
$ws.Range("C5").Value = 11796.5722254491
$ws.Range("E5").Value = 8343.96739787434
$ws.Range("F5").Value = 336.028317638475

$ws.Range("B6").Value = 12534.66680153
$ws.Range("C6").Value = 11161.0767346955
$ws.Range("E6").Value = 8450.90890760041
$ws.Range("F6").Value = 314.005235095662

$ws.Range("B7").Value = 4831.69084056906
$ws.Range("C7").Value = 7800.19319521736
$ws.Range("E7").Value = 7908.83855191676
$ws.Range("F7").Value = 151.382156130588

$ws.Range("C9").Value = 11067.7874373768
$ws.Range("F9").Value = 380.514646050971

$ws.Range("C10").Value = 10666.7753333225
$ws.Range("F10").Value = 363.805808382042

$ws.Range("C11").Value = 10348.965861282
$ws.Range("F11").Value = 350.563747047023

$ws.Range("C12").Value = 10348.6432646394
$ws.Range("F12").Value = 350.550305520249

$ws.Range("C13").Value = 9891.60086855907
$ws.Range("F13").Value = 331.506872350233

$ws.Range("C14").Value = 7064.58628439619
$ws.Range("F14").Value = 198.174367717193

$ws.Range("C15").Value = 7100.03880974612
$ws.Range("F15").Value = 199.290581331905
